$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24
$ws.Range("O24").Value = '[''Northern Ireland'', ''Portugal'']'

# Row 54
$ws.Range("H54").Value = '[''Russia'', 3, -4, 2]'
$ws.Range("M54").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 5

# Row 55
$ws.Range("H55").Value = '[''Russia'', 3, -4, 2]'
$ws.Range("M55").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("Q55").Value = 5

# Row 56
$ws.Range("H56").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M56").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("Q56").Value = 5

# Row 57
$ws.Range("H57").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N57").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q57").Value = 6

# Row 58
$ws.Range("H58").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M58").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("Q58").Value = 7

# Row 59
$ws.Range("H59").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M59").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("Q59").Value = 7

# Row 60
$ws.Range("H60").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N60").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q60").Value = 8

# Row 61
$ws.Range("H61").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("M61").Value = '[''Switzerland'', ''Germany'', ''Ukraine'', ''Russia'']'
$ws.Range("Q61").Value = 9

# Row 62
$ws.Range("H62").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N62").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q62").Value = 10

# Row 63
$ws.Range("H63").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N63").Value = '[''Russia'', ''Spain'']'
$ws.Range("Q63").Value = 11

# Row 64
$ws.Range("H64").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N64").Value = '[''Ukraine'', ''Russia'']'
$ws.Range("Q64").Value = 12

# Row 65
$ws.Range("H65").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N65").Value = '[''Ukraine'', ''Russia'']'
$ws.Range("Q65").Value = 12

# Row 66
$ws.Range("H66").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N66").Value = '[''Ukraine'', ''Russia'']'
$ws.Range("Q66").Value = 12

# Row 67
$ws.Range("H67").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N67").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q67").Value = 13

# Row 68
$ws.Range("H68").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N68").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q68").Value = 13

# Row 69
$ws.Range("H69").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N69").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q69").Value = 13

# Row 70
$ws.Range("H70").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N70").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q70").Value = 13

# Row 71
$ws.Range("H71").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N71").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q71").Value = 13

# Row 72
$ws.Range("H72").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N72").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q72").Value = 13

# Row 73
$ws.Range("H73").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N73").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O73").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q73").Value = 13

# Row 74
$ws.Range("H74").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N74").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O74").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q74").Value = 13

# Row 75
$ws.Range("H75").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N75").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O75").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q75").Value = 13

# Row 76
$ws.Range("H76").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N76").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O76").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q76").Value = 13

# Row 77
$ws.Range("H77").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N77").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O77").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q77").Value = 14

# Row 78
$ws.Range("H78").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N78").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O78").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q78").Value = 14

# Row 79
$ws.Range("H79").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N79").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O79").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q79").Value = 14

# Row 80
$ws.Range("H80").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N80").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O80").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q80").Value = 15

# Row 81
$ws.Range("H81").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N81").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O81").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q81").Value = 16

# Row 82
$ws.Range("H82").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N82").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O82").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q82").Value = 17

# Row 83
$ws.Range("H83").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N83").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O83").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q83").Value = 18

# Row 84
$ws.Range("H84").Value = '[''Russia'', 3, -5, 2]'
$ws.Range("N84").Value = '[''Russia'', ''Slovakia'']'
$ws.Range("O84").Value = '[''Slovakia'', ''Russia'']'
$ws.Range("Q84").Value = 19
